$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking columns must stay text (totalRuns, totalBalls, total4s, total6s, sr)
$ws.Range("G2:K14").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

# Row 2
$ws.Range("A2").Value = " Abu Dhabi"
$ws.Range("B2").Value = " October 25 2020"
$ws.Range("C2").Value = "Royals won by 8 wickets (with 10 balls remaining)"
$ws.Range("D2").Value = "Mumbai Indians"
$ws.Range("E2").Value = "Rajasthan Royals"
$ws.Range("F2").Value = "Hardik Pandya "
$ws.Range("G2").Value = "60"
$ws.Range("H2").Value = "21"
$ws.Range("I2").Value = "2"
$ws.Range("J2").Value = "7"
$ws.Range("K2").Value = "285.71"

# Row 3
$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " October 28 2020"
$ws.Range("C3").Value = "Mumbai won by 5 wickets (with 5 balls remaining)"
$ws.Range("D3").Value = "Mumbai Indians"
$ws.Range("E3").Value = "Royal Challengers Bangalore"
$ws.Range("F3").Value = "Hardik Pandya "
$ws.Range("G3").Value = "17"
$ws.Range("H3").Value = "15"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "2"
$ws.Range("K3").Value = "113.33"

# Row 4
$ws.Range("A4").Value = " Dubai (DSC)"
$ws.Range("B4").Value = " November 05 2020"
$ws.Range("C4").Value = "Mumbai won by 57 runs"
$ws.Range("D4").Value = "Mumbai Indians"
$ws.Range("E4").Value = "Delhi Capitals"
$ws.Range("F4").Value = "Hardik Pandya "
$ws.Range("G4").Value = "37"
$ws.Range("H4").Value = "14"
$ws.Range("I4").Value = "0"
$ws.Range("J4").Value = "5"
$ws.Range("K4").Value = "264.28"

# Row 5
$ws.Range("A5").Value = " Dubai (DSC)"
$ws.Range("B5").Value = " November 10 2020"
$ws.Range("C5").Value = "Mumbai won by 5 wickets (with 8 balls remaining)"
$ws.Range("D5").Value = "Mumbai Indians"
$ws.Range("E5").Value = "Delhi Capitals"
$ws.Range("F5").Value = "Hardik Pandya "
$ws.Range("G5").Value = "3"
$ws.Range("H5").Value = "5"
$ws.Range("I5").Value = "0"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "60.00"

# Row 6
$ws.Range("A6").Value = " Abu Dhabi"
$ws.Range("B6").Value = " October 16 2020"
$ws.Range("C6").Value = "Mumbai won by 8 wickets (with 19 balls remaining)"
$ws.Range("D6").Value = "Mumbai Indians"
$ws.Range("E6").Value = "Kolkata Knight Riders"
$ws.Range("F6").Value = "Hardik Pandya "
$ws.Range("G6").Value = "21"
$ws.Range("H6").Value = "11"
$ws.Range("I6").Value = "3"
$ws.Range("J6").Value = "1"
$ws.Range("K6").Value = "190.90"

# Row 7
$ws.Range("A7").Value = " Abu Dhabi"
$ws.Range("B7").Value = " October 01 2020"
$ws.Range("C7").Value = "Mumbai won by 48 runs"
$ws.Range("D7").Value = "Mumbai Indians"
$ws.Range("E7").Value = "Kings XI Punjab"
$ws.Range("F7").Value = "Hardik Pandya "
$ws.Range("G7").Value = "30"
$ws.Range("H7").Value = "11"
$ws.Range("I7").Value = "3"
$ws.Range("J7").Value = "2"
$ws.Range("K7").Value = "272.72"

# Row 8
$ws.Range("A8").Value = " Abu Dhabi"
$ws.Range("B8").Value = " September 19 2020"
$ws.Range("C8").Value = "Super Kings won by 5 wickets (with 4 balls remaining)"
$ws.Range("D8").Value = "Mumbai Indians"
$ws.Range("E8").Value = "Chennai Super Kings"
$ws.Range("F8").Value = "Hardik Pandya "
$ws.Range("G8").Value = "14"
$ws.Range("H8").Value = "10"
$ws.Range("I8").Value = "0"
$ws.Range("J8").Value = "2"
$ws.Range("K8").Value = "140.00"

# Row 9
$ws.Range("A9").Value = " Abu Dhabi"
$ws.Range("B9").Value = " October 06 2020"
$ws.Range("C9").Value = "Mumbai won by 57 runs"
$ws.Range("D9").Value = "Mumbai Indians"
$ws.Range("E9").Value = "Rajasthan Royals"
$ws.Range("F9").Value = "Hardik Pandya "
$ws.Range("G9").Value = "30"
$ws.Range("H9").Value = "19"
$ws.Range("I9").Value = "2"
$ws.Range("J9").Value = "1"
$ws.Range("K9").Value = "157.89"

# Row 10
$ws.Range("A10").Value = " Dubai (DSC)"
$ws.Range("B10").Value = " October 18 2020"
$ws.Range("C10").Value = "Match tied (Kings XI won the one-over eliminator)"
$ws.Range("D10").Value = "Mumbai Indians"
$ws.Range("E10").Value = "Kings XI Punjab"
$ws.Range("F10").Value = "Hardik Pandya "
$ws.Range("G10").Value = "8"
$ws.Range("H10").Value = "4"
$ws.Range("I10").Value = "0"
$ws.Range("J10").Value = "1"
$ws.Range("K10").Value = "200.00"

# Row 11
$ws.Range("A11").Value = " Dubai (DSC)"
$ws.Range("B11").Value = " September 28 2020"
$ws.Range("C11").Value = "Match tied (RCB won the one-over eliminator)"
$ws.Range("D11").Value = "Mumbai Indians"
$ws.Range("E11").Value = "Royal Challengers Bangalore"
$ws.Range("F11").Value = "Hardik Pandya "
$ws.Range("G11").Value = "15"
$ws.Range("H11").Value = "13"
$ws.Range("I11").Value = "0"
$ws.Range("J11").Value = "1"
$ws.Range("K11").Value = "115.38"

# Row 12
$ws.Range("A12").Value = " Abu Dhabi"
$ws.Range("B12").Value = " September 23 2020"
$ws.Range("C12").Value = "Mumbai won by 49 runs"
$ws.Range("D12").Value = "Mumbai Indians"
$ws.Range("E12").Value = "Kolkata Knight Riders"
$ws.Range("F12").Value = "Hardik Pandya "
$ws.Range("G12").Value = "18"
$ws.Range("H12").Value = "13"
$ws.Range("I12").Value = "2"
$ws.Range("J12").Value = "1"
$ws.Range("K12").Value = "138.46"

# Row 13
$ws.Range("A13").Value = " Abu Dhabi"
$ws.Range("B13").Value = " October 11 2020"
$ws.Range("C13").Value = "Mumbai won by 5 wickets (with 2 balls remaining)"
$ws.Range("D13").Value = "Mumbai Indians"
$ws.Range("E13").Value = "Delhi Capitals"
$ws.Range("F13").Value = "Hardik Pandya "
$ws.Range("G13").Value = "0"
$ws.Range("H13").Value = "2"
$ws.Range("I13").Value = "0"
$ws.Range("J13").Value = "0"
$ws.Range("K13").Value = "0.00"

# Row 14
$ws.Range("A14").Value = " Sharjah"
$ws.Range("B14").Value = " October 04 2020"
$ws.Range("C14").Value = "Mumbai won by 34 runs"
$ws.Range("D14").Value = "Mumbai Indians"
$ws.Range("E14").Value = "Sunrisers Hyderabad"
$ws.Range("F14").Value = "Hardik Pandya "
$ws.Range("G14").Value = "28"
$ws.Range("H14").Value = "19"
$ws.Range("I14").Value = "2"
$ws.Range("J14").Value = "2"
$ws.Range("K14").Value = "147.36"
